$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.846.48"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.808.54"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'309.51"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4652"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").Value = "'0.3693"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").Value = "'0.07360"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'0.8721"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'20.44"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.835.20"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'5.354"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").Value = "'6.505"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "'0.07052"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "'91.20"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "'0.000008698"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'14.72"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "26.866.41"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'5.314"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "  -4.09%  "
$ws.Range("D24").Value = "2.029.71"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'1.904"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'151.58"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'2.141"
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("D29").Value = "'5.300"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "'115.84"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "'0.08905"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "'0.7519"
$ws.Range("E32").Value = "  -5.24%  "
$ws.Range("D33").Value = "'1.151"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.918"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.454"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.01962"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'0.05256"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'2.427"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").Value = "'2.922"
$ws.Range("D42").Value = "'0.5293"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "'0.1662"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "'8.435"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "'0.4931"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").Value = "'10.29"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D49").Value = "'1.671"
$ws.Range("D50").Value = "'102.91"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "'0.06281"
$ws.Range("E51").Value = "  -1.63%  "
